$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new "BP" (villain/henchmen tag) column before column D; this pushes
# the existing D:N columns (and their widths/styles) over to E:O automatically.
$ws.Columns("D").Insert()

# New header + column width (narrow, no bestFit - matches the other tag columns)
$ws.Range("D1").Value = "BP"
$ws.Range("D1").ColumnWidth = 2.1666666666

# Per-row BP tag values (henchmen group "build points" cost)
$bpValues = @{
    2 = 3;  3 = 3;  4 = 3;  5 = 3;  6 = 3;  7 = 3;  8 = 3;  9 = 3;  10 = 4;
    11 = 3; 12 = 3; 13 = 3; 14 = 3; 15 = 3; 16 = 3; 17 = 3; 18 = 3; 19 = 2;
    20 = 3; 21 = 3; 22 = 2; 23 = 3; 24 = 1; 25 = 3; 26 = 3; 27 = 3; 28 = 3;
    29 = 3; 30 = 3; 31 = 3; 32 = 3; 33 = 3; 34 = 3; 35 = 3; 36 = 3; 37 = 3;
    38 = 3; 39 = 3; 40 = 2; 41 = 3; 42 = 2
}
foreach ($row in $bpValues.Keys) {
    $ws.Range("D$row").Value = $bpValues[$row]
}

# Restore view/selection (frozen header pane scrolled back to top, B4 selected)
$ws.Range("B4").Select() | Out-Null
